$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new "2022-Q1" worksheet right after "2021-Q4" (i.e. right
#    before "总计"). Duplicate the "2021-Q4" sheet so the new sheet
#    inherits the same sheet-level properties (outline/page setup,
#    margins, header style, column A numbering style, etc.) and then
#    overwrite its data with the 2022-Q1 fund holdings.
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2021-Q4")
$refSheet.Copy($null, $refSheet)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The copied sheet has 6 data rows (A2:A7); the new data only needs 4
# (A2:A5), so drop the extra two rows before filling in values.
$newSheet.Range("6:7").Delete()

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# The source data stores fund code / scale / position figures as plain text
# (not numbers) -- format these columns as text first so values like the
# zero-padded fund code "003318" survive the assignment untouched.
$newSheet.Range("B2:G5").NumberFormat = "@"

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Range("B2").Value = "003318"
$newSheet.Range("C2").Value = "景顺长城中证500行业中性低波动指数"
$newSheet.Range("D2").Value = "13.99"
$newSheet.Range("E2").Value = "93.88"
$newSheet.Range("F2").Value = "1.21"
$newSheet.Range("G2").Value = "0.1693"
$newSheet.Range("H2").Value = 6

$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Range("B3").Value = "005357"
$newSheet.Range("C3").Value = "富国国企改革灵活配置混合"
$newSheet.Range("D3").Value = "1.13"
$newSheet.Range("E3").Value = "87.21"
$newSheet.Range("F3").Value = "2.41"
$newSheet.Range("G3").Value = "0.0272"
$newSheet.Range("H3").Value = 7

$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Range("B4").Value = "160639"
$newSheet.Range("C4").Value = "鹏华中证高铁产业指数（LOF）"
$newSheet.Range("D4").Value = "0.89"
$newSheet.Range("E4").Value = "94.72"
$newSheet.Range("F4").Value = "2.40"
$newSheet.Range("G4").Value = "0.0214"
$newSheet.Range("H4").Value = 9

$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Range("B5").Value = "512260"
$newSheet.Range("C5").Value = "华安中证500行业中性低波动ETF"
$newSheet.Range("D5").Value = "1.17"
$newSheet.Range("E5").Value = "96.94"
$newSheet.Range("F5").Value = "1.25"
$newSheet.Range("G5").Value = "0.0146"
$newSheet.Range("H5").Value = 6

# Drop the temporary "Text" number format now that the values are locked in
# as text, so the cells end up with plain/default formatting (matching the
# unstyled data cells used elsewhere in the workbook).
$newSheet.Range("B2:G5").Style = "Normal"

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new top data row for
#    2022-Q1, shifting the existing quarterly rows down by one, and
#    renumber the index column (A).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.23

# Renumber the index column (A) for the rows that got pushed down.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(7, 1).Value = 5

# Restore the originally-active sheet/tab selection (the workbook opened
# on "2020-Q4"), since only the sheet data changed, not the active view.
$wb.Worksheets.Item("2020-Q4").Activate()
